# fix: rm duplicate item from timeline fig
# Remove the row containing the duplicate "Improve sampling from OMs" entry
# (month=11, year=2020) from the timeline table. All subsequent rows shift
# up by one and the now-unused shared string is dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 holds: 11 | 2020 | Improve sampling from OMs | planned
$ws.Rows.Item(5).Delete()

# Leave the selection where Excel would naturally land after this edit.
$ws.Range("C13").Select()
